$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Update column widths
# Note: Excel's ColumnWidth (character units) and the stored OOXML <col width>
# differ by a constant offset of 5/6 (~0.8333) with the default Calibri 11 font,
# so we subtract that offset to land exactly on the target stored width.
$colWidthOffset = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 11 - $colWidthOffset
$ws.Columns.Item(2).ColumnWidth = 11 - $colWidthOffset
$ws.Columns.Item(6).ColumnWidth = 32 - $colWidthOffset
$ws.Columns.Item(8).ColumnWidth = 61 - $colWidthOffset

# Update data rows 2-11
Set-TextCell $ws "A2" "Elmo"
Set-TextCell $ws "B2" "Tosi"
Set-TextCell $ws "C2" "24/05/2002"
Set-TextCell $ws "D2" "Millepini"
Set-TextCell $ws "E2" "PSCNNR73L22C712C"
Set-TextCell $ws "F2" "giancarlogolgi@example.org"
Set-TextCell $ws "G2" "0565876678"
Set-TextCell $ws "H2" "Incrocio Viridiana, 82 Appartamento 98`n08015, Macomer (NU)"
Set-TextCell $ws "A3" "Umberto"
Set-TextCell $ws "B3" "Veltroni"
Set-TextCell $ws "C3" "23/06/1993"
Set-TextCell $ws "D3" "Meri'"
Set-TextCell $ws "E3" "LFRMRT64R20I294N"
Set-TextCell $ws "F3" "giannuzzibenedetto@example.com"
Set-TextCell $ws "G3" "0789106263"
Set-TextCell $ws "H3" "Incrocio Gloria, 988 Appartamento 68`n36065, Mussolente (VI)"
Set-TextCell $ws "A4" "Enrico"
Set-TextCell $ws "B4" "Saragat"
Set-TextCell $ws "C4" "21/02/1992"
Set-TextCell $ws "D4" "Mottola"
Set-TextCell $ws "E4" "MLPSTN26S12E804I"
Set-TextCell $ws "F4" "rossettiflavia@example.com"
Set-TextCell $ws "G4" "377242544"
Set-TextCell $ws "H4" "Via Guariento, 77`n84014, Nocera Inferiore (SA)"
Set-TextCell $ws "A5" "Piero"
Set-TextCell $ws "B5" "Lerner"
Set-TextCell $ws "C5" "24/06/1996"
Set-TextCell $ws "D5" "Celzi"
Set-TextCell $ws "E5" "PCLSVN76P44D185J"
Set-TextCell $ws "F5" "nicola58@example.org"
Set-TextCell $ws "G5" "+39 0824513384"
Set-TextCell $ws "H5" "Contrada Alessia, 951 Piano 2`n13835, Botto (BI)"
Set-TextCell $ws "A6" "Dolores"
Set-TextCell $ws "B6" "Comolli"
Set-TextCell $ws "C6" "26/09/1999"
Set-TextCell $ws "D6" "Acquacanina"
Set-TextCell $ws "E6" "GNNVTR03B56H634C"
Set-TextCell $ws "F6" "graziellaroth@example.com"
Set-TextCell $ws "G6" "0371482077"
Set-TextCell $ws "H6" "Vicolo Federico, 28 Piano 3`n20035, Villa Cortese (MI)"
Set-TextCell $ws "A7" "Armando"
Set-TextCell $ws "B7" "Gibilisco"
Set-TextCell $ws "C7" "11/12/1996"
Set-TextCell $ws "D7" "San Pietro Val Lemina"
Set-TextCell $ws "E7" "ZCHPNI99L31A766E"
Set-TextCell $ws "F7" "guglielmosoderini@example.org"
Set-TextCell $ws "G7" "0572513888"
Set-TextCell $ws "H7" "Vicolo Danilo, 26 Piano 4`n74021, Carosino (TA)"
Set-TextCell $ws "A8" "Galasso"
Set-TextCell $ws "B8" "Travaglia"
Set-TextCell $ws "C8" "06/05/2005"
Set-TextCell $ws "D8" "Vaglie"
Set-TextCell $ws "E8" "MRCTMT74A42B371M"
Set-TextCell $ws "F8" "carmelo04@example.net"
Set-TextCell $ws "G8" "0942019544"
Set-TextCell $ws "H8" "Incrocio Bianca, 29`n28857, Santa Maria Maggiore (VB)"
Set-TextCell $ws "A9" "Durante"
Set-TextCell $ws "B9" "Brugnaro"
Set-TextCell $ws "C9" "04/07/2004"
Set-TextCell $ws "D9" "Valestra"
Set-TextCell $ws "E9" "PDRSRN42H19L406I"
Set-TextCell $ws "F9" "coriolanoovadia@example.net"
Set-TextCell $ws "G9" "076557727"
Set-TextCell $ws "H9" "Vicolo Ansaldo, 12`n58037, Santa Fiora (GR)"
Set-TextCell $ws "A10" "Flavia"
Set-TextCell $ws "B10" "Bonatti"
Set-TextCell $ws "C10" "25/10/2002"
Set-TextCell $ws "D10" "Case Di Nava"
Set-TextCell $ws "E10" "FRNMRL92A68A373V"
Set-TextCell $ws "F10" "lboldu@example.net"
Set-TextCell $ws "G10" "+39 057342862"
Set-TextCell $ws "H10" "Viale Baldassare, 98`n00079, Colle Di Fuori (RM)"
Set-TextCell $ws "A11" "Raffaello"
Set-TextCell $ws "B11" "Gigli"
Set-TextCell $ws "C11" "17/02/2005"
Set-TextCell $ws "D11" "San Lorenzo Pioppa"
Set-TextCell $ws "E11" "BRRDNI61C18L453Z"
Set-TextCell $ws "F11" "doriaclaudio@example.net"
Set-TextCell $ws "G11" "+39 351556726"
Set-TextCell $ws "H11" "Canale Interminelli, 7`n18021, Borgomaro (IM)"

# Cells in column H contain embedded line breaks, which makes the engine
# mark the row height as an explicit/custom value. Re-running AutoFit
# brings the rows back to the (non-custom) default-computed height so the
# saved XML doesn't carry a stray ht/customHeight attribute.
$ws.Range("2:11").EntireRow.AutoFit()
